$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the existing G4 timestamp (recorded a bit later than before).
$ws.Range("G4").Value = 42606.498657407406

# Append three new trade rows (5, 6, 7) with the same shape as the existing data.
$rows = @(
    @{ Row = 5; A = 9820.09;             B = 9875.39; C = 316.81; D = 318.58999999999997; E = $true; F = 0.56000000000000005; G = 42606.585590277777; H = $false },
    @{ Row = 6; A = 9752.33;             B = 9820.09; C = 316.81; D = 319;                 E = $true; F = 0.69;                 G = 42606.586701388886; H = $false },
    @{ Row = 7; A = 9697.7199999999993;  B = 9752.33; C = 316.81; D = 318.58999999999997; E = $true; F = 0.56000000000000005; G = 42606.587812500002; H = $false }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value = $r.A
    $ws.Cells.Item($n, 2).Value = $r.B
    $ws.Cells.Item($n, 3).Value = $r.C
    $ws.Cells.Item($n, 4).Value = $r.D
    $ws.Cells.Item($n, 5).Value = $r.E
    $ws.Cells.Item($n, 6).Value = $r.F
    $ws.Cells.Item($n, 7).Value = $r.G
    $ws.Cells.Item($n, 8).Value = $r.H
}
